# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45332 (2024-02-10) to 45333 (2024-02-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 27 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45332) {
        $cell.Value = 45333
    }
}
